$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.627.43'
$ws.Range("E2").Value = '  +0.48%  '
$ws.Range("D3").Value = '1.852.42'
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.027'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.41%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.91'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.026'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4379'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3800'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07408'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8833'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.08%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.57'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.87%  '
$ws.Range("D12").Value = '1.856.75'
$ws.Range("E12").Value = '  +0.66%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.525'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.722'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07145'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '85.14'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.032'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009084'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.87%  '
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.52%  '
$ws.Range("D21").Value = '27.661.86'
$ws.Range("E21").Value = '  +0.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.289'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.25'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.33%  '
$ws.Range("D24").Value = '2.097.12'
$ws.Range("E24").Value = '  +1.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.012'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.43'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.34%  '
$ws.Range("E27").Value = '  +0.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.345'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.985'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.49'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08993'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7766'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.93%  '
$ws.Range("E33").Value = '  +0.84%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.991'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.559'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.73%  '
$ws.Range("E36").Value = '  -0.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.140'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01974'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05274'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.48%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.862'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5190'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.53%  '
$ws.Range("E42").Value = '  +1.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.855'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.851'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '110.14'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.08%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.71'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.06604'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.029'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.704'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4709'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.898'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.01%  '
